# Collapse the blank paragraphs that surround two of the standalone
# text paragraphs: "The cat is sleeping on the wardrobe" loses the
# blank paragraph before and after it, and the blank paragraph between
# the "Weather is cloudy outside" paragraph and the next text
# paragraph is removed as well.

$d = $word.ActiveDocument

# Walk the paragraph collection looking for a paragraph whose own text
# is empty but whose next sibling is the text we care about, then pull
# that text up into the (currently empty) paragraph's run and drop the
# now-redundant paragraphs.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    # Strip the trailing paragraph mark before comparing.
    $text = $text.Substring(0, $text.Length - 1)

    if ($text -eq "") {
        $next = $para.Next()
        if ($next -ne $null) {
            $nextText = $next.Range.Text
            $nextText = $nextText.Substring(0, $nextText.Length - 1)

            if ($nextText -eq "The cat is sleeping on the wardrobe") {
                # Move the text into this (blank) paragraph's run, then
                # delete the paragraph that used to hold it and the
                # blank paragraph that trailed it.
                $para.Range.Text = $nextText
                $next.Range.Delete()
                $after = $d.Paragraphs.Item($i + 1)
                $afterText = $after.Range.Text
                $afterText = $afterText.Substring(0, $afterText.Length - 1)
                if ($afterText -eq "") {
                    $after.Range.Delete()
                }
                break
            }
        }
    }
}

# Now remove the blank paragraph that immediately follows the
# "Weather is cloudy outside" paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $text = $text.Substring(0, $text.Length - 1)

    if ($text -eq "Weather is cloudy outside") {
        $next = $para.Next()
        if ($next -ne $null) {
            $nextText = $next.Range.Text
            $nextText = $nextText.Substring(0, $nextText.Length - 1)
            if ($nextText -eq "") {
                $next.Range.Delete()
            }
        }
        break
    }
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [$($p.Range.Text)]"
}
